# Regenerate column G ("K") values on Sheet1 (rows 2-64) to reflect the
# freshly recalculated strikeout counts (K) that replace the previous
# "Strike#" derived figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for G2:G64, in row order (row 2 first).
$newK = @(
    1,1,0,2,2,2,1,1,1,0,1,1,2,0,0,2,0,2,1,1,
    0,1,0,1,0,0,3,0,1,0,2,0,1,2,1,0,2,5,0,0,
    1,2,2,1,4,1,1,3,2,3,2,0,2,1,0,0,2,1,0,2,
    1,2,1
)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
